# docs: updated existing Architecture, Logic, Storage, UI component Class Diagrams
#
# Logic Component Class Diagram (single-slide deck):
#   - "AddressBook" label (two-line label "AddressBook" / "Parser") renamed
#     to "ModsUni" (label becomes "ModsUni" / "Parser").
#   - The "XYZCommand = AddCommand, FindCommand, etc." note updates its
#     examples: AddCommand -> LoginCommand and FindCommand -> SaveCommand.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the two shapes we need to touch by scanning their text, so the
# script isn't brittle against shape re-ordering.
$addressBookShape = $null
$noteShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }
    $t = $shp.TextFrame.TextRange.Text
    if ($t -eq "AddressBookParser") {
        $addressBookShape = $shp
    } elseif ($t -like "*AddCommand*" -and $t -like "*FindCommand*") {
        $noteShape = $shp
    }
}
if ($addressBookShape -eq $null) { $addressBookShape = $s.Shapes.Item(9) }
if ($noteShape -eq $null) { $noteShape = $s.Shapes.Item(48) }

# --- "AddressBook" / "Parser" label -> "ModsUni" / "Parser" ---------------
$abRange = $addressBookShape.TextFrame.TextRange
$abPara1 = $abRange.Paragraphs(1)
$abPara1.Text = "ModsUni"

# --- "XYZCommand = AddCommand, FindCommand, etc." note --------------------
$noteRange = $noteShape.TextFrame.TextRange

$noteText = $noteRange.Text
$addIdx = $noteText.IndexOf("AddCommand")
$noteRange.Characters($addIdx + 1, "AddCommand".Length).Text = "LoginCommand"

$noteText = $noteRange.Text
$findIdx = $noteText.IndexOf("FindCommand")
$noteRange.Characters($findIdx + 1, "FindCommand".Length).Text = "SaveCommand"
